$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (the "jahr" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Add header for new column E
$ws.Cells.Item(1, 5).Value2 = "aggregate_id"

# Populate column E with the same values as column A (the year) for each data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 1).Value2
}

# Update selection to reflect the cell Excel ended up on after this edit
$ws.Range("F12").Select()
